$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): strip leading apostrophes, reword the "Increase" header ---
$ws.Range("A1").Value = "Property Name"
$ws.Range("B1").Value = "Property Value"
$ws.Range("C1").Value = "Increase value?"
# D1 ("Remarks") and F1 (instructions) are unchanged.

# --- Column A (property names): strip the leading apostrophes ---
$ws.Range("A2").Value  = "AUTHOR_NAME"
$ws.Range("A3").Value  = "Angle_Units"
$ws.Range("A4").Value  = "DN"
$ws.Range("A5").Value  = "DRAWING_TITLE"
$ws.Range("A6").Value  = "FC-DATE"
$ws.Range("A7").Value  = "FC-REV"
$ws.Range("A8").Value  = "FC-SC"
$ws.Range("A9").Value  = "FC-SH"
$ws.Range("A10").Value = "FC-SI"
$ws.Range("A11").Value = "FreeCAD_DRAWING"
$ws.Range("A12").Value = "Length_Units"
$ws.Range("A13").Value = "Mass_Units"
$ws.Range("A14").Value = "Number of sheets"
$ws.Range("A15").Value = "PN"
$ws.Range("A16").Value = "SI-1"

# --- Column B: drop the unit-marker values, keep the apostrophe-free label ---
$ws.Range("B3").ClearContents()
$ws.Range("B11").Value = "FreeCAD DRAWING"
$ws.Range("B12").ClearContents()
$ws.Range("B13").ClearContents()

# --- Column C: clear every "No"/"Yes" value - the Yes/No list is gone ---
$ws.Range("C2:C16").ClearContents()
$ws.Range("C9").Value = "X"

# Remove the old data validation (list "No, Yes") on C2:C16
$ws.Range("C2:C16").Validation.Delete()

# C1's header cell keeps the "quote prefix" style that the workbook already
# used elsewhere (e.g. B5) - copy formats only from B5 onto C1.
$ws.Range("B5").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Selection moves to B13 ---
$ws.Range("B13").Select()

# --- Best-effort: shrink/re-position the workbook window to match the author's view ---
try {
    $win = $wb.Windows.Item(1)
    $win.Left = -120
    $win.Top = -120
    $win.Width = 29040
    $win.Height = 15720
} catch {
    # Window geometry is a cosmetic, host-level setting; ignore if unsupported.
}
